$p = $ppt.ActivePresentation
$s = $p.Slides.Item(4)
$shp = $s.Shapes.Item(2)
$tr = $shp.TextFrame.TextRange

# Paragraph 6: "Cloud Platform: Microsoft Azure" -> merged Platform + Core Service
$tr.Paragraphs(6, 1).Runs(1, 1).Text = "Platform: Microsoft Azure with Sentinel (SIEM/SOAR)"

# Paragraph 7: "Core Service: Microsoft Sentinel (SIEM/SOAR)" -> merged Data Store + Analytics
$tr.Paragraphs(7, 1).Runs(1, 1).Text = "Data & Analytics: Log Analytics Workspace with KQL queries"

# Paragraph 8 ("Data Store: Log Analytics Workspace") is now redundant/merged away - remove it
$tr.Paragraphs(8, 1).Delete()

# Remaining paragraph (now index 8, previously 9): tighten the Integration wording
$tr.Paragraphs(8, 1).Runs(1, 1).Text = "Integration: Microsoft Defender, Azure Monitor, third-party tools"
